$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -----------------
# This shared string is used by Overview!E2, Overview!F2, zh-cn!C2 and
# de-de!C2. Update every cell that carries it so the runtime collapses them
# back onto a single (new) shared string, mirroring the <si> text edit in
# the diff rather than forking a separate string for just one cell.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ---------------------------------------------------
# Overview columns E & F (the zh-cn / de-de status columns) and the "Status"
# column (C) on both the zh-cn and de-de detail sheets all shrink from
# ~17.22 stored width down to ~13.41 stored width.
#
# This runtime quantizes stored column width to steps of 1/6, so feed it the
# ColumnWidth (character-unit) value whose rounded result lands on the
# nearest achievable step to the target stored width of 13.4101845877511
# (i.e. 13.333333333333334 == 80/6); 12.5 sits safely in the middle of the
# input range that rounds to that step.
$targetColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth       # column C
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth       # column C
